$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 5 new rows (150-154) to the "anotacoes" table, describing lesson 63
# "MockEmailService com Logger. Padroes Strategy e Template Method" (section
# 5 - Serviço de email). This mirrors the pattern used by the previous block
# of rows (146-149, lesson 62).
# ---------------------------------------------------------------------------

# Copy the formatting (styles, number formats, wrap text, etc.) of the last
# existing data row down onto the five new rows so they look identical to
# the rest of the table.
$ws.Range("B149:G149").Copy() | Out-Null
$ws.Range("B150:G154").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# nome aula (column E) shared by all 5 new rows
$nomeAula = "MockEmailService com Logger. Padroes Strategy e Template Method"

# Fill in B/C/D (Seção, Nome da Seção, Aula) and F (abordagem da aula) for
# all five rows first, then E (nome aula) last, so new shared-string entries
# get appended in the same order Excel produced them in originally.

# Row 150 --------------------------------------------------------------
$ws.Cells.Item(150, 2).Value = 5
$ws.Cells.Item(150, 3).Value = "Serviço de email"
$ws.Cells.Item(150, 4).Value = 63
$ws.Cells.Item(150, 6).Value = "1:57`n5. Serviço de email`n63. MockEmailService com Logger. Padroes Strategy e Template Method`npadrão de projeto `"Strategy`" - polimorfismo"
$ws.Cells.Item(150, 7).Value = "`n`n`n`n`n`n`n`n"

# Row 151 --------------------------------------------------------------
$ws.Cells.Item(151, 2).Value = 5
$ws.Cells.Item(151, 3).Value = "Serviço de email"
$ws.Cells.Item(151, 4).Value = 63
$ws.Cells.Item(151, 6).Value = "3:48`n5. Serviço de email`n63. MockEmailService com Logger. Padroes Strategy e Template Method`npadrao de projeto `"Template Method`""

# Row 152 --------------------------------------------------------------
$ws.Cells.Item(152, 2).Value = 5
$ws.Cells.Item(152, 3).Value = "Serviço de email"
$ws.Cells.Item(152, 4).Value = 63
$ws.Cells.Item(152, 6).Value = "9:23`n5. Serviço de email`n63. MockEmailService com Logger. Padroes Strategy e Template Method`ncriação do MockMailService - email de mentirinha no log do servidor"

# Row 153 --------------------------------------------------------------
$ws.Cells.Item(153, 2).Value = 5
$ws.Cells.Item(153, 3).Value = "Serviço de email"
$ws.Cells.Item(153, 4).Value = 63
$ws.Cells.Item(153, 6).Value = "9:56`n5. Serviço de email`n63. MockEmailService com Logger. Padroes Strategy e Template Method`ninstancia um objeto de Logger"

# Row 154 --------------------------------------------------------------
$ws.Cells.Item(154, 2).Value = 5
$ws.Cells.Item(154, 3).Value = "Serviço de email"
$ws.Cells.Item(154, 4).Value = 63
$ws.Cells.Item(154, 6).Value = "12:29`n5. Serviço de email`n63. MockEmailService com Logger. Padroes Strategy e Template Method`ncriação de metodo @Bean de test na classe TestConfig - retornando um MockemailService"

# nome aula (column E), added last so it becomes the final new shared string
$ws.Cells.Item(150, 5).Value = $nomeAula
$ws.Cells.Item(151, 5).Value = $nomeAula
$ws.Cells.Item(152, 5).Value = $nomeAula
$ws.Cells.Item(153, 5).Value = $nomeAula
$ws.Cells.Item(154, 5).Value = $nomeAula

# Row heights, matching how Excel auto-sized these wrapped cells.
$ws.Rows.Item(150).RowHeight = 135
$ws.Rows.Item(151).RowHeight = 75
$ws.Rows.Item(152).RowHeight = 90
$ws.Rows.Item(153).RowHeight = 75
$ws.Rows.Item(154).RowHeight = 90

# ---------------------------------------------------------------------------
# Grow the Excel table ("Tabela1") and its AutoFilter to cover the new rows.
# ---------------------------------------------------------------------------
$table = $ws.ListObjects.Item("Tabela1")
$table.Resize($ws.Range("B1:G154"))

# ---------------------------------------------------------------------------
# Update the sheet view so the new last row is visible/selected, matching
# what Excel does after entering data at the bottom of the sheet.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 149
$ws.Range("E153").Select()
